# Weights.xlsx update - "More work on JPADCore_v2 concerning weights analysis."
# Updates the computed weight-estimation values on each results sheet.

$wb = $excel.ActiveWorkbook

# --- GLOBAL RESULTS --------------------------------------------------------
$ws = $wb.Worksheets.Item("GLOBAL RESULTS")
$ws.Range("C2").Value  = 3999.999999999999
$ws.Range("C6").Value  = 27512.318491806276
$ws.Range("C7").Value  = 269803.6781376719
$ws.Range("C8").Value  = 23673.118491806283
$ws.Range("C9").Value  = 232153.98745767202
$ws.Range("C10").Value = 21662.206642625653
$ws.Range("C11").Value = 212433.67877190482
$ws.Range("C18").Value = 19114.626711865414
$ws.Range("C19").Value = 187450.4540439149
$ws.Range("C20").Value = 18718.626711865414
$ws.Range("C21").Value = 183567.02064391493
$ws.Range("C22").Value = 11986.626711865418
$ws.Range("C23").Value = 117548.65284391496
$ws.Range("C24").Value = 10757.083066365416
$ws.Range("C25").Value = 105490.94865277238
$ws.Range("C26").Value = 10739.123066365419
$ws.Range("C27").Value = 105314.82121877241
$ws.Range("C28").Value = 1017.96
$ws.Range("C29").Value = 9982.777433999996
$ws.Range("C30").Value = 1000.0
$ws.Range("C31").Value = 9806.649999999998
$ws.Range("C32").Value = 11986.626711865418
$ws.Range("C33").Value = 117548.65284391496

# --- FUSELAGE ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("FUSELAGE")
$ws.Range("C6").Value  = 1863.0
$ws.Range("D6").Value  = -44.23157516613781
$ws.Range("C7").Value  = 1876.0
$ws.Range("D7").Value  = -43.84242351673351
$ws.Range("C8").Value  = 3217.0
$ws.Range("D8").Value  = -3.6999341435670083
$ws.Range("C9").Value  = 1987.0
$ws.Range("D9").Value  = -40.51966712566604
$ws.Range("C12").Value = 1863.0

# --- WING ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("WING")
$ws.Range("C2").Value  = 22000.0
$ws.Range("C7").Value  = 2454.0
$ws.Range("D7").Value  = -88.84545454545454
$ws.Range("C8").Value  = 1921.0
$ws.Range("D8").Value  = -91.26818181818182
$ws.Range("C9").Value  = 1641.0
$ws.Range("D9").Value  = -92.5409090909091
$ws.Range("D10").Value = -89.55
$ws.Range("C11").Value = 2413.0
$ws.Range("D11").Value = -89.03181818181818
$ws.Range("C12").Value = 2271.0
$ws.Range("D12").Value = -89.67727272727272
$ws.Range("C13").Value = 1641.0

# --- HORIZONTAL TAIL ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("HORIZONTAL TAIL")
$ws.Range("C2").Value  = 3299.9999999999995
$ws.Range("D7").Value  = -92.18181818181817
$ws.Range("C8").Value  = 150.0
$ws.Range("D8").Value  = -95.45454545454544
$ws.Range("C9").Value  = 144.0
$ws.Range("D9").Value  = -95.63636363636363
$ws.Range("C10").Value = 258.0

# --- VERTICAL TAIL ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("VERTICAL TAIL")
$ws.Range("C2").Value = 3299.9999999999995
$ws.Range("D7").Value = -89.99999999999999
$ws.Range("C8").Value = 232.0
$ws.Range("D8").Value = -92.96969696969695
$ws.Range("C9").Value = 330.0

# --- NACELLES ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("NACELLES")
$ws.Range("C3").Value  = 712.0
$ws.Range("D3").Value  = -13.043478260869545
$ws.Range("C10").Value = 356.0
$ws.Range("D10").Value = -13.043478260869561
$ws.Range("C12").Value = 356.0
$ws.Range("C17").Value = 356.0
$ws.Range("D17").Value = -13.043478260869561
$ws.Range("C19").Value = 356.0

# --- POWER PLANT ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("POWER PLANT")
$ws.Range("D3").Value = -13.043478260869545

# --- LANDING GEARS ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("LANDING GEARS")
$ws.Range("C3").Value  = 1057.4352685771673
$ws.Range("C6").Value  = 911.0
$ws.Range("D6").Value  = 34.80319621189702
$ws.Range("C7").Value  = 1135.0
$ws.Range("D7").Value  = 67.94909736608466
$ws.Range("C8").Value  = 1281.0
$ws.Range("D8").Value  = 89.55312222551052
$ws.Range("C9").Value  = 1057.0
$ws.Range("D9").Value  = 56.407221071322894
$ws.Range("C10").Value = 911.0
